# Lab3 ("лр3") block: append 16 new rows (57-72) to Лист1 with the same
# layout/pattern as the existing лр1 / лр2 blocks, mirroring the author's
# "remove indexing from filenames" change that introduced a third lab.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$values = @(8545, 1288, 2614, 7170, 2709, 3402, 8781, 8771, 3091, 2195, 1184, 2667, 5917, 3929, 3083, 7546)

$firstRow = 57
$lastRow = 72

# Column B values for the new rows.
for ($i = 0; $i -lt $values.Length; $i++) {
    $row = $firstRow + $i
    $ws.Cells.Item($row, 2).Value = $values[$i]
}

# Column C: same relative formula as the лр2 block; Excel will regroup this
# into its own shared-formula block when written back out.
$ws.Range("C$firstRow`:C$lastRow").Formula = '=IF(B57=$E$1,1,0)'

# Column A: label for the new block + matching style (centered, like A2 and
# A17) copied from the existing лр1 block so no new cell-format entries are
# introduced, then merge the whole block like the other two labs.
$ws.Range("A2").Copy() | Out-Null
$ws.Range("A$firstRow`:A$lastRow").PasteSpecial(-4122) | Out-Null
$ws.Range("A$firstRow").Value = "лр3"
$ws.Range("A$firstRow`:A$lastRow").Merge() | Out-Null

# Match the author's final view state (scrolled/selected near the new rows).
$win = $excel.ActiveWindow
$win.SplitRow = 1
$win.ScrollRow = 43
$ws.Range("F53").Select() | Out-Null
